$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "convnext_tiny"
$ws.Cells.Item(8, 3).Value = 0.9523809523809523
$ws.Cells.Item(8, 5).Value = 0.975609756097561
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(10, 4).Value = 0.95
$ws.Cells.Item(10, 5).Value = 0.9743589743589743
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(16, 3).Value = 0.9523809523809523
$ws.Cells.Item(16, 5).Value = 0.975609756097561
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(23, 1).Value = "densenet121"
$ws.Cells.Item(24, 4).Value = 0.95
$ws.Cells.Item(24, 5).Value = 0.9743589743589743
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(27, 5).Value = 1
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).Value = 0.95
$ws.Cells.Item(29, 5).Value = 0.9743589743589743
$ws.Cells.Item(30, 3).Value = 0.9523809523809523
$ws.Cells.Item(30, 5).Value = 0.975609756097561
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(34, 5).Value = 1
$ws.Cells.Item(35, 3).Value = 0.9090909090909091
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(35, 5).Value = 0.9523809523809523
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(37, 5).Value = 1
$ws.Cells.Item(40, 3).Value = 0.9523809523809523
$ws.Cells.Item(40, 5).Value = 0.975609756097561
$ws.Cells.Item(41, 4).Value = 0.9
$ws.Cells.Item(41, 5).Value = 0.9473684210526315
$ws.Cells.Item(42, 3).Value = 1
$ws.Cells.Item(42, 5).Value = 1
$ws.Cells.Item(44, 1).Value = "efficientnet_b0"
$ws.Cells.Item(48, 3).Value = 0.95
$ws.Cells.Item(48, 4).Value = 0.95
$ws.Cells.Item(48, 5).Value = 0.95
$ws.Cells.Item(50, 3).Value = 0.9523809523809523
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(50, 5).Value = 0.975609756097561
$ws.Cells.Item(51, 3).Value = 1
$ws.Cells.Item(51, 5).Value = 1
$ws.Cells.Item(52, 4).Value = 1
$ws.Cells.Item(52, 5).Value = 1
$ws.Cells.Item(56, 3).Value = 1
$ws.Cells.Item(56, 5).Value = 1
$ws.Cells.Item(58, 3).Value = 1
$ws.Cells.Item(58, 5).Value = 1
$ws.Cells.Item(65, 1).Value = "efficientnet_b3"
$ws.Cells.Item(69, 3).Value = 0.9523809523809523
$ws.Cells.Item(69, 4).Value = 1
$ws.Cells.Item(69, 5).Value = 0.975609756097561
$ws.Cells.Item(71, 3).Value = 0.9523809523809523
$ws.Cells.Item(71, 4).Value = 1
$ws.Cells.Item(71, 5).Value = 0.975609756097561
$ws.Cells.Item(77, 3).Value = 1
$ws.Cells.Item(77, 4).Value = 0.95
$ws.Cells.Item(77, 5).Value = 0.9743589743589743
$ws.Cells.Item(83, 4).Value = 0.95
$ws.Cells.Item(83, 5).Value = 0.9743589743589743
$ws.Cells.Item(86, 1).Value = "resnet101"
$ws.Cells.Item(87, 4).Value = 0.95
$ws.Cells.Item(87, 5).Value = 0.9743589743589743
$ws.Cells.Item(92, 4).Value = 0.95
$ws.Cells.Item(92, 5).Value = 0.9743589743589743
$ws.Cells.Item(93, 3).Value = 1
$ws.Cells.Item(93, 5).Value = 1
$ws.Cells.Item(95, 4).Value = 0.95
$ws.Cells.Item(95, 5).Value = 0.9743589743589743
$ws.Cells.Item(98, 3).Value = 0.9523809523809523
$ws.Cells.Item(98, 5).Value = 0.975609756097561
$ws.Cells.Item(102, 3).Value = 0.9523809523809523
$ws.Cells.Item(102, 5).Value = 0.975609756097561
$ws.Cells.Item(103, 3).Value = 0.9523809523809523
$ws.Cells.Item(103, 5).Value = 0.975609756097561
$ws.Cells.Item(104, 4).Value = 1
$ws.Cells.Item(104, 5).Value = 1
$ws.Cells.Item(107, 1).Value = "resnet50"
$ws.Cells.Item(109, 3).Value = 0.9523809523809523
$ws.Cells.Item(109, 5).Value = 0.975609756097561
$ws.Cells.Item(113, 3).Value = 1
$ws.Cells.Item(113, 4).Value = 0.95
$ws.Cells.Item(113, 5).Value = 0.9743589743589743
$ws.Cells.Item(116, 4).Value = 0.95
$ws.Cells.Item(116, 5).Value = 0.9743589743589743
$ws.Cells.Item(119, 3).Value = 0.9523809523809523
$ws.Cells.Item(119, 4).Value = 1
$ws.Cells.Item(119, 5).Value = 0.975609756097561
$ws.Cells.Item(125, 4).Value = 1
$ws.Cells.Item(125, 5).Value = 1
$ws.Cells.Item(127, 3).Value = 1
$ws.Cells.Item(127, 5).Value = 1
$ws.Cells.Item(134, 4).Value = 1
$ws.Cells.Item(134, 5).Value = 1
$ws.Cells.Item(136, 4).Value = 0.95
$ws.Cells.Item(136, 5).Value = 0.9743589743589743
$ws.Cells.Item(140, 3).Value = 1
$ws.Cells.Item(140, 5).Value = 1
$ws.Cells.Item(142, 3).Value = 0.9523809523809523
$ws.Cells.Item(142, 5).Value = 0.975609756097561
$ws.Cells.Item(146, 3).Value = 1
$ws.Cells.Item(146, 4).Value = 1
$ws.Cells.Item(146, 5).Value = 1
$ws.Cells.Item(147, 4).Value = 1
$ws.Cells.Item(147, 5).Value = 1
$ws.Cells.Item(148, 3).Value = 1
$ws.Cells.Item(148, 5).Value = 1
$ws.Cells.Item(149, 1).Value = "vit_b_16"
$ws.Cells.Item(153, 3).Value = 0.9473684210526315
$ws.Cells.Item(153, 4).Value = 0.9
$ws.Cells.Item(153, 5).Value = 0.9230769230769231
$ws.Cells.Item(155, 3).Value = 1
$ws.Cells.Item(155, 4).Value = 0.85
$ws.Cells.Item(155, 5).Value = 0.918918918918919
$ws.Cells.Item(160, 3).Value = 0.9523809523809523
$ws.Cells.Item(160, 5).Value = 0.975609756097561
$ws.Cells.Item(161, 3).Value = 0.8695652173913043
$ws.Cells.Item(161, 4).Value = 1
$ws.Cells.Item(161, 5).Value = 0.9302325581395349
$ws.Cells.Item(163, 4).Value = 0.95
$ws.Cells.Item(163, 5).Value = 0.9743589743589743
$ws.Cells.Item(167, 4).Value = 0.95
$ws.Cells.Item(167, 5).Value = 0.9743589743589743
$ws.Cells.Item(168, 3).Value = 0.9090909090909091
$ws.Cells.Item(168, 5).Value = 0.9523809523809523
$ws.Cells.Item(169, 3).Value = 0.9523809523809523
$ws.Cells.Item(169, 5).Value = 0.975609756097561
